# "numbers suck but trying to fix them"
# Birat BF sheet: insert a new "fueltype" column before the existing
# %-fuel-mix columns (J..O), label it, fill it with "PCI coal" for the
# fuel-mix rows, and correct the BF-bb fuel-mix numbers in row 10
# (B/C/D/G/H/I) to the revised PCI-coal-based figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Birat BF")

# Insert a new column at J (10th column); everything from the old J
# onward (fossil/biofuel % columns, secondary-fuel-type, etc.) shifts
# right by one, from J..O to K..P.
$ws.Columns.Item(10).Insert()

# New column header + sub-header.
$ws.Range("J1").Value2 = "fueltype"
$ws.Range("J2").Value2 = "t BF fuel/t fresh steel"

# Every fuel-mix data row (4 through 10) uses PCI coal in the new column.
$ws.Range("J4:J10").Value2 = "PCI coal"

# Row 10 (the "BF bb" scenario row) gets corrected fuel-mix figures.
$ws.Range("B10").Value2 = 0
$ws.Range("C10").Formula = "=0.1169+0.0731"
$ws.Range("D10").Formula = "=0.305/(1-C10)"
$ws.Range("G10").Formula = "=0"
$ws.Range("H10").Formula = "=0.478342726/(1-C10)"
$ws.Range("I10").Formula = "=0.1376816/(1-C10)"
